$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.624.20"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "'1.634.08"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'212.32"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "'0.520"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'23.33"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").Value = "'0.0870"
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "'1.867.33"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "'1.640.03"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "'0.553"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "'65.25"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "'27.609.17"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'230.83"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "'0.0₃0720"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "'7.58"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'10.63"
$ws.Range("E22").Value = "  +5.77%  "
$ws.Range("D23").Value = "'4.35"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("E24").Value = "  +7.62%  "
$ws.Range("D25").Value = "'149.64"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "'6.89"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "'0.0484"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Value = "'1.477.63"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").Value = "'1.55"
$ws.Range("E35").Value = "  -1.34%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").Value = "'0.940"
$ws.Range("E37").Value = "  +4.79%  "
$ws.Range("D38").Value = "'0.880"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "'0.0167"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'68.03"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").Value = "'2.20"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("E46").Value = "  -4.40%  "
$ws.Range("D47").Value = "'1.775.09"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "'1.75"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").Value = "'87.62"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").Value = "'0.0992"
$ws.Range("E51").Value = "  +0.88%  "
